$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under the duplicate_image_filename column (column E)
# for every data row in the first table (rows 2-21).
foreach ($r in 2..21) {
    $ws.Range("E$r").Value = "NA"
}

# Workaround: touching the sheet can cause the previously-blank F1 cell
# (an empty string placeholder) to be re-serialized with a stray value on
# save. Explicitly clear it so it round-trips back to being blank.
$ws.Range("F1").ClearContents()
